$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (stored as a serial date number, advance by one day)
$ws.Range("A1").Value = 45311

# Update prices in the price list
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 223.526
